# Updates the cryptos list table (columns B:E, rows 2-50) to reflect refreshed
# market data, as captured by the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, never letting Excel reinterpret
# numeric-looking strings (e.g. "7.190", "0.00000000123") as numbers, which
# would silently drop trailing zeros / switch to scientific notation.
# NumberFormat "@" forces text entry, then ClearFormats() removes the
# temporary text-format style again so the cell keeps its original (default)
# style, matching the source workbook exactly.
function Set-CellText($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-CellText 2 4 '29.169.74'

Set-CellText 3 4 '1.840.19'
Set-CellText 3 5 '  +0.62%  '

Set-CellText 4 5 '  +0.50%  '

Set-CellText 5 4 '244.39'
Set-CellText 5 5 '  +0.28%  '

Set-CellText 6 4 '0.6187'
Set-CellText 6 5 '  -2.10%  '

Set-CellText 7 4 '1.005'
Set-CellText 7 5 '  +0.61%  '

Set-CellText 8 4 '0.07495'
Set-CellText 8 5 '  -0.23%  '

Set-CellText 9 4 '0.2930'
Set-CellText 9 5 '  -0.28%  '

Set-CellText 10 4 '23.23'
Set-CellText 10 5 '  +1.45%  '

Set-CellText 11 4 '0.07714'
Set-CellText 11 5 '  -0.09%  '

Set-CellText 12 4 '1.845.36'
Set-CellText 12 5 '  +0.66%  '

Set-CellText 13 4 '5.011'
Set-CellText 13 5 '  +0.39%  '

Set-CellText 14 4 '0.6734'
Set-CellText 14 5 '  +0.44%  '

Set-CellText 15 4 '82.93'
Set-CellText 15 5 '  -0.18%  '

Set-CellText 16 4 '0.000009306'
Set-CellText 16 5 '  -4.08%  '

Set-CellText 17 4 '5.955'
Set-CellText 17 5 '  -1.97%  '

Set-CellText 18 4 '29.174.35'
Set-CellText 18 5 '  +0.53%  '

Set-CellText 19 4 '2.121.60'
Set-CellText 19 5 '  +1.76%  '

Set-CellText 20 4 '232.55'
Set-CellText 20 5 '  +2.68%  '

Set-CellText 21 4 '12.68'

Set-CellText 22 5 '  +0.74%  '

Set-CellText 23 4 '7.190'
Set-CellText 23 5 '  +0.18%  '

Set-CellText 24 4 '1.005'
Set-CellText 24 5 '  +0.65%  '

Set-CellText 25 4 '160.45'
Set-CellText 25 5 '  +0.50%  '

Set-CellText 26 4 '8.545'
Set-CellText 26 5 '  +0.10%  '

Set-CellText 27 4 '0.1392'

Set-CellText 28 4 '17.87'
Set-CellText 28 5 '  -0.15%  '

Set-CellText 29 4 '1.508'
Set-CellText 29 5 '  +0.62%  '

Set-CellText 30 4 '4.179'
Set-CellText 30 5 '  +1.48%  '

Set-CellText 31 4 '4.141'
Set-CellText 31 5 '  +1.70%  '

Set-CellText 32 4 '0.05565'
Set-CellText 32 5 '  +3.64%  '

Set-CellText 33 4 '1.204'
Set-CellText 33 5 '  +0.18%  '

Set-CellText 34 4 '0.7501'
Set-CellText 34 5 '  +0.97%  '

Set-CellText 35 4 '1.846'
Set-CellText 35 5 '  -0.68%  '

Set-CellText 36 5 '  +0.41%  '

Set-CellText 37 4 '2.668'
Set-CellText 37 5 '  +0.56%  '

Set-CellText 38 4 '2.777'
Set-CellText 38 5 '  +1.04%  '

Set-CellText 39 4 '1.229.36'
Set-CellText 39 5 '  -1.25%  '

Set-CellText 40 4 '0.01786'
Set-CellText 40 5 '  -0.03%  '

Set-CellText 41 4 '6.501'
Set-CellText 41 5 '  -1.35%  '

Set-CellText 42 4 '0.8979'
Set-CellText 42 5 '  -0.71%  '

Set-CellText 43 5 '  +0.58%  '

Set-CellText 44 4 '2.017.97'
Set-CellText 44 5 '  +1.67%  '

Set-CellText 45 4 '102.11'
Set-CellText 45 5 '  +0.53%  '

Set-CellText 46 4 '66.02'
Set-CellText 46 5 '  +1.83%  '

Set-CellText 47 2 'BabyDogeCoin'
Set-CellText 47 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText 47 4 '0.00000000123'
Set-CellText 47 5 '  +0.06%  '

Set-CellText 48 2 'Mantle'
Set-CellText 48 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-CellText 48 4 '0.5109'
Set-CellText 48 5 '  +0.24%  '

Set-CellText 49 4 '0.4085'
Set-CellText 49 5 '  +0.42%  '

Set-CellText 50 4 '9.132'
Set-CellText 50 5 '  +1.37%  '
